# The "Tower of Hanoi:" title slide and the "Tower of Hanoi interface:"
# detail slide were out of order in the deck (the detail slide was placed
# before its own title slide). Fix the order by moving the title slide
# ("Tower of Hanoi:", currently slide 9) so it comes right before the
# interface/detail slide (currently slide 8).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$s.MoveTo(8)
